# daily auto push: 2026-02-24 10:04 UTC
# A new measurement row is inserted above row 865 (2026/02/24, weekday 火,
# value 17, ranking 201). All the existing rows from 865 through 906 shift
# down by one (to 866..907), and the sheet's used range grows from
# A1:D906 to A1:D907.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push everything from row 865 down by one row, leaving a blank row 865.
$ws.Rows.Item(865).Insert()

# The "日付" column stores dates as plain text (e.g. "2026/02/24"), not as
# real Excel date serials, elsewhere in this sheet. Force text formatting
# before writing so Excel doesn't auto-convert the string into a date
# value/format, then drop back to the sheet's normal (unstyled) look so the
# new row matches its neighbours.
$ws.Range("A865").NumberFormat = "@"
$ws.Range("A865").Value = "2026/02/24"
$ws.Range("A865").Style = "Normal"

$ws.Range("B865").Value = "火"
$ws.Range("C865").Value = 17
$ws.Range("D865").Value = 201
